# Weekly fruit/vegetable price update: a new observation is inserted as
# row 30 ("Rabanito" @ Vega Modelo de Temuco), pushing the existing rows
# 30..59 down to 31..60 (dimension grows from A1:R59 to A1:R60).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 30 - shifts rows 30..59 down to
# 31..60 and carries the D-column (date) number formatting along.
$ws.Rows.Item(30).Insert()

# Populate the newly-inserted row 30 with the new weekly record.
$ws.Range("A30").Value = 10
$ws.Range("B30").Value = "Vega Modelo de Temuco"
$ws.Range("C30").Value = "La Araucanía"
$ws.Range("D30").Value = 44638
$ws.Range("E30").Value = 9
$ws.Range("F30").Value = 300000001
$ws.Range("G30").Value = "Rabanito"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 20
$ws.Range("K30").Value = 7000
$ws.Range("L30").Value = 7000
$ws.Range("M30").Value = 7000
$ws.Range("N30").Value = "$/docena de paquetes"
$ws.Range("O30").Value = "Provincia de Cautín"
$ws.Range("P30").Value = 583
$ws.Range("Q30").Value = 12
$ws.Range("R30").Value = "Hortaliza"
